# Weekly fruit/vegetable update: insert a new weekly price record for
# "Feria Lagunitas de Puerto Montt - Betarraga" at row 226, pushing the
# existing rows 226:240 down to 227:241.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 226 (existing row 226 and below shift down to 227+)
$ws.Rows.Item(226).Insert()

# Populate the newly inserted row with the new weekly data point.
$ws.Cells.Item(226, 1).Value  = 4
$ws.Cells.Item(226, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(226, 3).Value  = "Los Lagos"
$ws.Cells.Item(226, 4).Value  = 44585
$ws.Cells.Item(226, 5).Value  = 10
$ws.Cells.Item(226, 6).Value  = 100114014
$ws.Cells.Item(226, 7).Value  = "Betarraga"
$ws.Cells.Item(226, 8).Value  = "Sin especificar"
$ws.Cells.Item(226, 9).Value  = "Primera"
$ws.Cells.Item(226, 10).Value = 250
$ws.Cells.Item(226, 11).Value = 1000
$ws.Cells.Item(226, 12).Value = 1000
$ws.Cells.Item(226, 13).Value = 1000
$ws.Cells.Item(226, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(226, 15).Value = "Región del Maule"
$ws.Cells.Item(226, 16).Value = 200
$ws.Cells.Item(226, 17).Value = 5
$ws.Cells.Item(226, 18).Value = "Hortaliza"
